# Weekly driver report update for 2025-05-05
# Updates the "Good Drivers" table (rows 13-24) on the active sheet with
# refreshed counts from the latest pull, which also re-orders a few rows
# (14/15/16 rotate, 18/19 swap, 20/21 swap) to match the new sort produced
# by the upstream report generator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-update) contents for rows 13-24, columns A-J.
# Columns: A=adapter-driver, B=good sum, C=critical sum, D=warning sum,
#          E=client count, F=total sum, G=adapter, H=driver,
#          I=good roaming calc (%), J=driver vintage
$rows = @(
    @{ Row=13; A="Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.0.4"; B=1293197; C=4322; D=1009; E=1990; F=1298528; H="22.250.0.4"; J="2023-07-25" },
    @{ Row=14; A="Intel(R) Wi-Fi 6E AX211 160MHz - 22.220.0.4"; B=31517;   C=112;  D=0;    E=66;   F=31629;   H="22.220.0.4"; J="2023-03-28" },
    @{ Row=15; A="Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8";  B=467311; C=772;  D=931;  E=706;  F=469014;  H="23.10.0.8";  J="2023-10-30" },
    @{ Row=16; A="Intel(R) Wi-Fi 6E AX211 160MHz - 23.120.0.3"; B=455081; C=1861; D=52;   E=639;  F=456994;  H="23.120.0.3"; J="2025-02-05" },
    @{ Row=17; A="Intel(R) Wi-Fi 6E AX211 160MHz - 22.230.0.8"; B=1787924;C=3326; D=2614; E=3038; F=1793864; H="22.230.0.8"; J="2023-05-08" },
    @{ Row=18; A="Intel(R) Wi-Fi 6E AX211 160MHz - 23.70.2.3";  B=218767; C=334;  D=313;  E=573;  F=219414;  H="23.70.2.3";  J="2024-07-23" },
    @{ Row=19; A="Intel(R) Wi-Fi 6E AX211 160MHz - 22.110.1.1"; B=135467; C=189;  D=263;  E=196;  F=135919;  H="22.110.1.1"; J="2022-01-01" },
    @{ Row=20; A="Intel(R) Wi-Fi 6E AX211 160MHz - 23.100.0.4"; B=240434; C=421;  D=37;   E=409;  F=240892;  H="23.100.0.4"; J="2024-11-10" },
    @{ Row=21; A="Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.1.3";  B=151287; C=285;  D=75;   E=332;  F=151647;  H="23.80.1.3";  J="2024-09-03" },
    @{ Row=22; A="Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"; B=272039; C=213;  D=131;  E=316;  F=272383;  H="22.100.1.1"; J="2022-05-01" },
    @{ Row=23; A="Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"; B=14561;  C=0;    D=0;    E=59;   F=14561;   H="22.150.0.3"; J="2022-05-23" },
    @{ Row=24; A="Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"; B=12018;  C=0;    D=0;    E=61;   F=12018;   H="22.150.3.1"; J="2022-08-29" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    # Keep the "driver vintage" column as plain text (matches source data,
    # which stores dates as inline strings, not real dates).
    $ws.Cells.Item($rowNum, 10).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 10).Value = $r.J
}
